$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Perfil de Solicitação de Empréstimo / Condição de Pagador
$ws.Range("V2").Value = "Baixo"
$ws.Range("W2").Value = "Médio Pagador"

# Row 3
$ws.Range("V3").Value = "Alto"
$ws.Range("W3").Value = "Mal Pagador"

# Row 4 (only Condição de Pagador changes)
$ws.Range("W4").Value = "Médio Pagador"

# Row 5
$ws.Range("V5").Value = "Baixo"
$ws.Range("W5").Value = "Médio Pagador"

# Row 6 (only Perfil changes)
$ws.Range("V6").Value = "Sem Risco"

# Row 7 (only Perfil changes)
$ws.Range("V7").Value = "Médio"

# Row 9 (only Perfil changes)
$ws.Range("V9").Value = "Alto"

# Row 10 (only Perfil changes)
$ws.Range("V10").Value = "Médio"

# Row 11 (only Perfil changes)
$ws.Range("V11").Value = "Sem Risco"
